# Update the "datetimeFigureOut" date placeholder text from 9/8/21 to
# 9/23/21 across the slide master and every slide layout (11 layouts),
# matching the commit's "Updated examples for unioffice v1.16.0" edit.

$p = $ppt.ActivePresentation
$newDate = "9/23/21"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
        }
        if ($isDatePlaceholder) {
            if ($sh.TextFrame.TextRange.Text -eq "9/8/21") {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster

# Slide master's own Date Placeholder shape.
Update-DatePlaceholder $master.Shapes

# Every custom layout (slideLayout1.xml .. slideLayout11.xml) owned by
# the slide master has its own Date Placeholder shape too.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
